$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header values for new columns P1 (14) and Q1 (15)
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# Copy the header style (bold/centered/bordered) from O1 to P1:Q1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update data rows 2-25, columns B..Q
# Row 2
$ws.Range("B2").Value = 13.21888947808733
$ws.Range("C2").Value = 9.148687272497256
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 28.21366595912454
$ws.Range("F2").Value = 30.27884324296193
$ws.Range("G2").Value = 28.40185174476629
$ws.Range("H2").Value = 2.458333511100583
$ws.Range("I2").Value = 2.618761695987909
$ws.Range("J2").Value = 10.2699439176619
$ws.Range("K2").Value = 15.45449320215411
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 12.23748714003568
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 12.68728345213712
$ws.Range("Q2").Value = 0

# Row 3
$ws.Range("B3").Value = 12.34559179291352
$ws.Range("C3").Value = 8.564469868069075
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 26.34472502459244
$ws.Range("F3").Value = 28.16301984326747
$ws.Range("G3").Value = 27.70669695370652
$ws.Range("H3").Value = 2.22056518867392
$ws.Range("I3").Value = 2.779447459452341
$ws.Range("J3").Value = 10.2331900850231
$ws.Range("K3").Value = 15.43526853197001
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 11.48066808198653
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 12.83935181526205
$ws.Range("Q3").Value = 0

# Row 4
$ws.Range("B4").Value = 11.77463857204742
$ws.Range("C4").Value = 8.189473222262366
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 25.13088557130743
$ws.Range("F4").Value = 26.79314218379242
$ws.Range("G4").Value = 27.27969715696966
$ws.Range("H4").Value = 2.069524846383745
$ws.Range("I4").Value = 2.882221785106803
$ws.Range("J4").Value = 10.21288083802414
$ws.Range("K4").Value = 15.42622293951482
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 10.98951906652294
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 12.93375331481854
$ws.Range("Q4").Value = 0

# Row 5
$ws.Range("B5").Value = 11.52685733960785
$ws.Range("C5").Value = 8.039602459871965
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 24.61924732691535
$ws.Range("F5").Value = 26.21694803853475
$ws.Range("G5").Value = 27.08558523298714
$ws.Range("H5").Value = 2.006485484619084
$ws.Range("I5").Value = 2.928193691965745
$ws.Range("J5").Value = 10.20166513319635
$ws.Range("K5").Value = 15.41569558429189
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 10.78118399214159
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 12.97156015873806
$ws.Range("Q5").Value = 0

# Row 6
$ws.Range("B6").Value = 11.4774336948624
$ws.Range("C6").Value = 8.023075950458873
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 24.5329700040801
$ws.Range("F6").Value = 26.1201820470107
$ws.Range("G6").Value = 27.02866674807851
$ws.Range("H6").Value = 1.995698227946642
$ws.Range("I6").Value = 2.939808305046902
$ws.Range("J6").Value = 10.19557280088175
$ws.Range("K6").Value = 15.40482181045799
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 10.74428882191887
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 12.9767188806658
$ws.Range("Q6").Value = 0

# Row 7
$ws.Range("B7").Value = 11.75067428392937
$ws.Range("C7").Value = 8.210426002245466
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 25.12322767100302
$ws.Range("F7").Value = 26.78544435617556
$ws.Range("G7").Value = 27.21018227266019
$ws.Range("H7").Value = 2.068031445828518
$ws.Range("I7").Value = 2.89323037435845
$ws.Range("J7").Value = 10.20111812545347
$ws.Range("K7").Value = 15.40108362750074
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 10.98165174866688
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 12.93122624231862
$ws.Range("Q7").Value = 0

# Row 8
$ws.Range("B8").Value = 12.89972972834286
$ws.Range("C8").Value = 8.97871273879782
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 27.581940600936
$ws.Range("F8").Value = 29.56389055882486
$ws.Range("G8").Value = 28.07646995337808
$ws.Range("H8").Value = 2.376569134632406
$ws.Range("I8").Value = 2.686753934139764
$ws.Range("J8").Value = 10.24159161564938
$ws.Range("K8").Value = 15.41442478198947
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 11.97576940267674
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 12.73580107151512
$ws.Range("Q8").Value = 0

# Row 9
$ws.Range("B9").Value = 14.91205247318888
$ws.Range("C9").Value = 10.31977262603857
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 31.88209804503983
$ws.Range("F9").Value = 34.4587799414611
$ws.Range("G9").Value = 29.88334761591804
$ws.Range("H9").Value = 2.94365952348296
$ws.Range("I9").Value = 2.654286852139454
$ws.Range("J9").Value = 10.36033197679083
$ws.Range("K9").Value = 15.50842975318631
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 13.72320640496621
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 12.36492631994635
$ws.Range("Q9").Value = 0

# Row 10
$ws.Range("B10").Value = 16.17516792113765
$ws.Range("C10").Value = 11.16008572334624
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 33.84714229353244
$ws.Range("F10").Value = 37.72874519847217
$ws.Range("G10").Value = 30.84497388078908
$ws.Range("H10").Value = 3.301878126092816
$ws.Range("I10").Value = 2.90978997973153
$ws.Range("J10").Value = 10.39625749455128
$ws.Range("K10").Value = 15.48774801193306
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 14.82382282033805
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 12.08902408739418
$ws.Range("Q10").Value = 0

# Row 11
$ws.Range("B11").Value = 16.32738545922534
$ws.Range("C11").Value = 10.90964334685155
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 27.18207987237584
$ws.Range("F11").Value = 39.14805722232856
$ws.Range("G11").Value = 28.39393715490085
$ws.Range("H11").Value = 3.766635742312189
$ws.Range("I11").Value = 2.955714248827072
$ws.Range("J11").Value = 9.8962306622363
$ws.Range("K11").Value = 14.64388430997038
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 14.84014283225323
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 11.94166569559194
$ws.Range("Q11").Value = 0

# Row 12
$ws.Range("B12").Value = 16.22586323287459
$ws.Range("C12").Value = 10.50906537760554
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 21.13288116032132
$ws.Range("F12").Value = 39.67592721363786
$ws.Range("G12").Value = 26.20384585552203
$ws.Range("H12").Value = 4.724221533570519
$ws.Range("I12").Value = 2.944999591855198
$ws.Range("J12").Value = 9.48810376336871
$ws.Range("K12").Value = 13.99065926443797
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 14.6403023730495
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 11.91362547864012
$ws.Range("Q12").Value = 0

# Row 13
$ws.Range("B13").Value = 15.89672788813648
$ws.Range("C13").Value = 9.985026108550251
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 15.01569306199732
$ws.Range("F13").Value = 39.56266403260141
$ws.Range("G13").Value = 23.96251328804266
$ws.Range("H13").Value = 5.87009141776567
$ws.Range("I13").Value = 2.894411325933132
$ws.Range("J13").Value = 9.109618063846925
$ws.Range("K13").Value = 13.4108246332328
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 14.25121273202812
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 11.96625216382721
$ws.Range("Q13").Value = 0

# Row 14
$ws.Range("B14").Value = 15.55475773141073
$ws.Range("C14").Value = 9.563474682454684
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 10.6664388171801
$ws.Range("F14").Value = 39.19167637503955
$ws.Range("G14").Value = 22.35564869558961
$ws.Range("H14").Value = 6.738672937385491
$ws.Range("I14").Value = 2.843059173640333
$ws.Range("J14").Value = 8.859166721067846
$ws.Range("K14").Value = 13.0413056265312
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 13.88828284636274
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 12.04124650054118
$ws.Range("Q14").Value = 0

# Row 15
$ws.Range("B15").Value = 15.41464383622755
$ws.Range("C15").Value = 9.43773783489092
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 9.587158293016538
$ws.Range("F15").Value = 38.96319309154519
$ws.Range("G15").Value = 21.93488233483884
$ws.Range("H15").Value = 6.938194496106381
$ws.Range("I15").Value = 2.822567207953838
$ws.Range("J15").Value = 8.801167810608517
$ws.Range("K15").Value = 12.9589777369205
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 13.75568452372774
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 12.07220937471136
$ws.Range("Q15").Value = 0

# Row 16
$ws.Range("B16").Value = 14.93048848041325
$ws.Range("C16").Value = 9.168866310127859
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 9.481679952207594
$ws.Range("F16").Value = 37.6346354302881
$ws.Range("G16").Value = 21.88022630983605
$ws.Range("H16").Value = 6.666071677465835
$ws.Range("I16").Value = 2.727617568961096
$ws.Range("J16").Value = 8.85766356009233
$ws.Range("K16").Value = 13.07768988196288
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 13.35341607387164
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 12.16466211712203
$ws.Range("Q16").Value = 0

# Row 17
$ws.Range("B17").Value = 14.74096478952206
$ws.Range("C17").Value = 9.193939817045486
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 11.81453417636688
$ws.Range("F17").Value = 36.80228778268388
$ws.Range("G17").Value = 22.69889964186271
$ws.Range("H17").Value = 5.944487749350684
$ws.Range("I17").Value = 2.681437709997931
$ws.Range("J17").Value = 9.031806542359137
$ws.Range("K17").Value = 13.35877577619024
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 13.23998720388946
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 12.19650878490917
$ws.Range("Q17").Value = 0

# Row 18
$ws.Range("B18").Value = 14.80764437891443
$ws.Range("C18").Value = 9.472755780563888
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 16.64003224081024
$ws.Range("F18").Value = 36.31710943239624
$ws.Range("G18").Value = 24.42243271128525
$ws.Range("H18").Value = 4.820466979548488
$ws.Range("I18").Value = 2.671960264912179
$ws.Range("J18").Value = 9.34052443819915
$ws.Range("K18").Value = 13.84358221907054
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 13.37590748848619
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 12.19439266303654
$ws.Range("Q18").Value = 0

# Row 19
$ws.Range("B19").Value = 15.04730753328827
$ws.Range("C19").Value = 9.976030702184174
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 23.07979269189842
$ws.Range("F19").Value = 36.15172494157628
$ws.Range("G19").Value = 26.66286682113545
$ws.Range("H19").Value = 3.682147010767213
$ws.Range("I19").Value = 2.707820991873096
$ws.Range("J19").Value = 9.727629078488349
$ws.Range("K19").Value = 14.44572160282087
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 13.6987699586383
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 12.1882074326891
$ws.Range("Q19").Value = 0

# Row 20
$ws.Range("B20").Value = 15.80242229812412
$ws.Range("C20").Value = 10.99582978913544
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 33.30127430104402
$ws.Range("F20").Value = 36.89155697731405
$ws.Range("G20").Value = 30.39322439660403
$ws.Range("H20").Value = 3.204845119628204
$ws.Range("I20").Value = 2.85282091038796
$ws.Range("J20").Value = 10.34876172025951
$ws.Range("K20").Value = 15.4125780499915
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 14.5280178739784
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 12.15400266679869
$ws.Range("Q20").Value = 0

# Row 21
$ws.Range("B21").Value = 16.80078600193536
$ws.Range("C21").Value = 11.72038314282004
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 36.08539221588159
$ws.Range("F21").Value = 39.30090300329764
$ws.Range("G21").Value = 31.67739939666336
$ws.Range("H21").Value = 3.530172555805767
$ws.Range("I21").Value = 3.0574715516422
$ws.Range("J21").Value = 10.47939928378353
$ws.Range("K21").Value = 15.56535958760306
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 15.41880498372669
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 11.95393542528966
$ws.Range("Q21").Value = 0

# Row 22
$ws.Range("B22").Value = 17.4206744089178
$ws.Range("C22").Value = 12.12005958453257
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 37.3944552075553
$ws.Range("F22").Value = 40.81974568026784
$ws.Range("G22").Value = 32.41839347270462
$ws.Range("H22").Value = 3.720871689818167
$ws.Range("I22").Value = 3.181497814652507
$ws.Range("J22").Value = 10.55178091510911
$ws.Range("K22").Value = 15.65198338252981
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 15.95529255480598
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 11.82328655499743
$ws.Range("Q22").Value = 0

# Row 23
$ws.Range("B23").Value = 17.11013264453392
$ws.Range("C23").Value = 11.88837490570064
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 36.70163037567293
$ws.Range("F23").Value = 40.01414225038423
$ws.Range("G23").Value = 32.09043498550345
$ws.Range("H23").Value = 3.620041136986892
$ws.Range("I23").Value = 3.111857417941077
$ws.Range("J23").Value = 10.52583310025378
$ws.Range("K23").Value = 15.63371141806441
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 15.6755363897325
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 11.8955265040085
$ws.Range("Q23").Value = 0

# Row 24
$ws.Range("B24").Value = 15.84895466083574
$ws.Range("C24").Value = 11.01304387250513
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 33.96529646859907
$ws.Range("F24").Value = 36.85121910575438
$ws.Range("G24").Value = 30.73927552396898
$ws.Range("H24").Value = 3.230059678456804
$ws.Range("I24").Value = 2.850101891996515
$ws.Range("J24").Value = 10.41326322305279
$ws.Range("K24").Value = 15.52856682236423
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 14.56320643759057
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 12.16587260250495
$ws.Range("Q24").Value = 0

# Row 25
$ws.Range("B25").Value = 14.36560102450955
$ws.Range("C25").Value = 10.00536226115283
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 30.77291795334995
$ws.Range("F25").Value = 33.19272459126713
$ws.Range("G25").Value = 29.28420607711359
$ws.Range("H25").Value = 2.793226313403603
$ws.Range("I25").Value = 2.560873182515035
$ws.Range("J25").Value = 10.3056237389514
$ws.Range("K25").Value = 15.4359118883556
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 13.26517399458283
$ws.Range("N25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 12.46003108938976
$ws.Range("Q25").Value = 0
